$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---

# Row 2: Execution Date changes
$ws.Cells.Item(2, 4).Value = "24/03/2022"

# Row 8: Execution Date + Test Result change
$ws.Cells.Item(8, 4).Value = "18/02/2022"
$ws.Cells.Item(8, 5).Value = "Pass"

# Row 10: Execution Date changes
$ws.Cells.Item(10, 4).Value = "17/02/2022"

# --- Append new rows ---

# Row 11
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "quote_ProposalSetup_227820_TC_03"
$ws.Cells.Item(11, 4).Value = "18/02/2022"
$ws.Cells.Item(11, 5).Value = "Fail"

# Row 12
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "quote_ProposalSetup_184671_TC_013"
$ws.Cells.Item(12, 4).Value = "22/02/2022"
$ws.Cells.Item(12, 5).Value = "Fail"

# Row 13
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "quote_ProposalSetup_227904_TC_012"
$ws.Cells.Item(13, 4).Value = "22/02/2022"
$ws.Cells.Item(13, 5).Value = "Pass"

# Row 14
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "quote_ProposalSetup_232246_TC_014"
$ws.Cells.Item(14, 4).Value = "25/03/2022"
$ws.Cells.Item(14, 5).Value = "Fail"

# Row 15
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "dealAnalysis_Finance_227824_TC_01"
$ws.Cells.Item(15, 4).Value = "31/03/2022"
$ws.Cells.Item(15, 5).Value = "Fail"
